$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values (recalculo de los resultados, cambios de LSPMW)
$ws.Range("B2").Value = 80448460319.36288
$ws.Range("C2").Value = 2412667512.985902
$ws.Range("D2").Value = 78035792806.37698

$ws.Range("B3").Value = 26997107502.8449
$ws.Range("C3").Value = 734245007.6737363
$ws.Range("D3").Value = 26262862495.17116

$ws.Range("B4").Value = 82260972169.52925
$ws.Range("C4").Value = 2476191769.954624
$ws.Range("D4").Value = 79784780399.57463

$ws.Range("B5").Value = 76853197667.61612
$ws.Range("C5").Value = 2267867024.890732
$ws.Range("D5").Value = 74585330642.72539

$ws.Range("B6").Value = 24167913345.41976
$ws.Range("C6").Value = 684729526.2772753
$ws.Range("D6").Value = 23483183819.14249

$ws.Range("B7").Value = 3220626050.680063
$ws.Range("C7").Value = 88451370.87369135
$ws.Range("D7").Value = 3132174679.806372

$ws.Range("B8").Value = 16236655844.35864
$ws.Range("C8").Value = 448083662.4119719
$ws.Range("D8").Value = 15788572181.94667

$ws.Range("B9").Value = 50817480527.65192
$ws.Range("C9").Value = 1410829212.825469
$ws.Range("D9").Value = 49406651314.82645

$ws.Range("B10").Value = 594.0351617115377
$ws.Range("C10").Value = 11.26879859863915
$ws.Range("D10").Value = 582.7663631128986
